$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 3016.5454  # ALC!H98: 3293.25 -> 3016.5454
$ws.Cells.Item(98, 9).Value = 3244.45  # ALC!I98: 3577.2222 -> 3244.45
$ws.Cells.Item(98, 11).Value = 3244.45  # ALC!K98: 3577.2222 -> 3244.45
$ws.Cells.Item(98, 13).Value = -1746.45  # ALC!M98: -2079.2222 -> -1746.45

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(116, 8).Value = 21649594  # ALC!H116: 22482192 -> 21649594
$ws.Cells.Item(116, 10).Value = 33337356  # ALC!J116: 37041276 -> 33337356
$ws.Cells.Item(116, 12).Value = 33337356  # ALC!L116: 37041276 -> 33337356
$ws.Cells.Item(116, 14).Value = -33344240  # ALC!N116: -37048160 -> -33344240

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(122, 8).Value = 3016.5454  # ALC!H122: 3293.25 -> 3016.5454
$ws.Cells.Item(122, 9).Value = 3244.45  # ALC!I122: 3577.2222 -> 3244.45
$ws.Cells.Item(122, 11).Value = 9733.349999999999  # ALC!K122: 10731.6666 -> 9733.349999999999
$ws.Cells.Item(122, 13).Value = -7283.349999999999  # ALC!M122: -8281.6666 -> -7283.349999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(133, 8).Value = 46936.25  # ALC!H133: 47436.25 -> 46936.25
$ws.Cells.Item(133, 10).Value = 46936.25  # ALC!J133: 47436.25 -> 46936.25
$ws.Cells.Item(133, 12).Value = 46936.25  # ALC!L133: 47436.25 -> 46936.25
$ws.Cells.Item(133, 14).Value = -57056.25  # ALC!N133: -57556.25 -> -57056.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 936.2  # ALC!H135: 998.36365 -> 936.2
$ws.Cells.Item(135, 9).Value = 730.8  # ALC!I135: 775 -> 730.8
$ws.Cells.Item(135, 11).Value = 6577.2  # ALC!K135: 6975 -> 6577.2
$ws.Cells.Item(135, 13).Value = -4042.2  # ALC!M135: -4440 -> -4042.2

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1931.1351  # ALC!H138: 1756.5714 -> 1931.1351
$ws.Cells.Item(138, 9).Value = 1644.0667  # ALC!I138: 1359.5 -> 1644.0667
$ws.Cells.Item(138, 10).Value = 2126.8635  # ALC!J138: 2117.5454 -> 2126.8635
$ws.Cells.Item(138, 11).Value = 4932.2001  # ALC!K138: 4078.5 -> 4932.2001
$ws.Cells.Item(138, 12).Value = 6380.5905  # ALC!L138: 6352.6362 -> 6380.5905
$ws.Cells.Item(138, 13).Value = 207.7999  # ALC!M138: 1061.5 -> 207.7999
$ws.Cells.Item(138, 14).Value = -16660.5905  # ALC!N138: -16632.6362 -> -16660.5905

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(139, 8).Value = 80814.836  # ALC!H139: 86397.8 -> 80814.836
$ws.Cells.Item(139, 10).Value = 92977.8  # ALC!J139: 102997.25 -> 92977.8
$ws.Cells.Item(139, 12).Value = 92977.8  # ALC!L139: 102997.25 -> 92977.8
$ws.Cells.Item(139, 14).Value = -103257.8  # ALC!N139: -113277.25 -> -103257.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 4349.75  # ARM!H10: 5376 -> 4349.75
$ws.Cells.Item(10, 9).Value = 4349.75  # ARM!I10: 5499.6665 -> 4349.75
$ws.Cells.Item(10, 10).Value = 0  # ARM!J10: 5005 -> 0
$ws.Cells.Item(10, 11).Value = 4349.75  # ARM!K10: 5499.6665 -> 4349.75
$ws.Cells.Item(10, 12).Value = 0  # ARM!L10: 5005 -> 0
$ws.Cells.Item(10, 13).Value = -4179.75  # ARM!M10: -5329.6665 -> -4179.75
$ws.Cells.Item(10, 14).ClearContents()  # ARM!N10: -5345 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6076.107  # ARM!H32: 5932.5386 -> 6076.107
$ws.Cells.Item(32, 9).Value = 4960.88  # ARM!I32: 5130.7085 -> 4960.88
$ws.Cells.Item(32, 10).Value = 15369.667  # ARM!J32: 15554.5 -> 15369.667
$ws.Cells.Item(32, 11).Value = 4960.88  # ARM!K32: 5130.7085 -> 4960.88
$ws.Cells.Item(32, 12).Value = 15369.667  # ARM!L32: 15554.5 -> 15369.667
$ws.Cells.Item(32, 13).Value = -4673.88  # ARM!M32: -4843.7085 -> -4673.88
$ws.Cells.Item(32, 14).Value = -15943.667  # ARM!N32: -16128.5 -> -15943.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4336.4287  # ARM!H61: 4523.846 -> 4336.4287
$ws.Cells.Item(61, 9).Value = 2185.7144  # ARM!I61: 2233.3333 -> 2185.7144
$ws.Cells.Item(61, 11).Value = 2185.7144  # ARM!K61: 2233.3333 -> 2185.7144
$ws.Cells.Item(61, 13).Value = -1973.7144  # ARM!M61: -2021.3333 -> -1973.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 58825856  # ARM!H74: 50002400 -> 58825856
$ws.Cells.Item(74, 9).Value = 100001050  # ARM!I74: 111112160 -> 100001050
$ws.Cells.Item(74, 10).Value = 4145.7144  # ARM!J74: 3509.6365 -> 4145.7144
$ws.Cells.Item(74, 11).Value = 100001050  # ARM!K74: 111112160 -> 100001050
$ws.Cells.Item(74, 12).Value = 4145.7144  # ARM!L74: 3509.6365 -> 4145.7144
$ws.Cells.Item(74, 13).Value = -100000176  # ARM!M74: -111111286 -> -100000176
$ws.Cells.Item(74, 14).Value = -5893.7144  # ARM!N74: -5257.636500000001 -> -5893.7144

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(77, 8).Value = 58825856  # ARM!H77: 50002400 -> 58825856
$ws.Cells.Item(77, 9).Value = 100001050  # ARM!I77: 111112160 -> 100001050
$ws.Cells.Item(77, 10).Value = 4145.7144  # ARM!J77: 3509.6365 -> 4145.7144
$ws.Cells.Item(77, 11).Value = 500005250  # ARM!K77: 555560800 -> 500005250
$ws.Cells.Item(77, 12).Value = 20728.572  # ARM!L77: 17548.1825 -> 20728.572
$ws.Cells.Item(77, 13).Value = -500000882  # ARM!M77: -555556432 -> -500000882
$ws.Cells.Item(77, 14).Value = -29464.572  # ARM!N77: -26284.1825 -> -29464.572

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 4188.3394  # ARM!H132: 4310.1665 -> 4188.3394
$ws.Cells.Item(132, 9).Value = 4169.522  # ARM!I132: 4318.1816 -> 4169.522
$ws.Cells.Item(132, 11).Value = 12508.566  # ARM!K132: 12954.5448 -> 12508.566
$ws.Cells.Item(132, 13).Value = -9978.565999999999  # ARM!M132: -10424.5448 -> -9978.565999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(133, 8).Value = 58816.332  # ARM!H133: 58849.668 -> 58816.332
$ws.Cells.Item(133, 10).Value = 58816.332  # ARM!J133: 58849.668 -> 58816.332
$ws.Cells.Item(133, 12).Value = 58816.332  # ARM!L133: 58849.668 -> 58816.332
$ws.Cells.Item(133, 14).Value = -63876.332  # ARM!N133: -63909.668 -> -63876.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4336.4287  # ARM!H136: 4523.846 -> 4336.4287
$ws.Cells.Item(136, 9).Value = 2185.7144  # ARM!I136: 2233.3333 -> 2185.7144
$ws.Cells.Item(136, 11).Value = 6557.1432  # ARM!K136: 6699.999899999999 -> 6557.1432
$ws.Cells.Item(136, 13).Value = -4007.1432  # ARM!M136: -4149.999899999999 -> -4007.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(137, 8).Value = 50778.5  # ARM!H137: 50518.668 -> 50778.5
$ws.Cells.Item(137, 9).Value = 0  # ARM!I137: 49999 -> 0
$ws.Cells.Item(137, 11).Value = 0  # ARM!K137: 49999 -> 0
$ws.Cells.Item(137, 13).ClearContents()  # ARM!M137: -44899 -> (removed)

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(138, 8).Value = 69381.5  # ARM!H138: 79688.25 -> 69381.5
$ws.Cells.Item(138, 10).Value = 89994  # ARM!J138: 89994.664 -> 89994
$ws.Cells.Item(138, 12).Value = 89994  # ARM!L138: 89994.664 -> 89994
$ws.Cells.Item(138, 14).Value = -100274  # ARM!N138: -100274.664 -> -100274

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(139, 8).Value = 55269.168  # ARM!H139: 56102.5 -> 55269.168
$ws.Cells.Item(139, 10).Value = 55269.168  # ARM!J139: 56102.5 -> 55269.168
$ws.Cells.Item(139, 12).Value = 55269.168  # ARM!L139: 56102.5 -> 55269.168
$ws.Cells.Item(139, 14).Value = -65549.16800000001  # ARM!N139: -66382.5 -> -65549.16800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(12, 8).Value = 1833.3334  # BSM!H12: 2650 -> 1833.3334
$ws.Cells.Item(12, 9).Value = 250  # BSM!I12: 300 -> 250
$ws.Cells.Item(12, 11).Value = 250  # BSM!K12: 300 -> 250
$ws.Cells.Item(12, 13).Value = -82  # BSM!M12: -132 -> -82

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 5000  # BSM!H63: 0 -> 5000
$ws.Cells.Item(63, 10).Value = 5000  # BSM!J63: 0 -> 5000
$ws.Cells.Item(63, 12).Value = 5000  # BSM!L63: 0 -> 5000
$ws.Cells.Item(63, 14).Value = -6372  # BSM!N63: None -> -6372

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(66, 8).Value = 5000  # BSM!H66: 0 -> 5000
$ws.Cells.Item(66, 10).Value = 5000  # BSM!J66: 0 -> 5000
$ws.Cells.Item(66, 12).Value = 15000  # BSM!L66: 0 -> 15000
$ws.Cells.Item(66, 14).Value = -21864  # BSM!N66: None -> -21864

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 9).Value = 35716524  # BSM!I94: 41669132 -> 35716524
$ws.Cells.Item(94, 10).Value = 3949.25  # BSM!J94: 3335 -> 3949.25
$ws.Cells.Item(94, 11).Value = 35716524  # BSM!K94: 41669132 -> 35716524
$ws.Cells.Item(94, 12).Value = 3949.25  # BSM!L94: 3335 -> 3949.25
$ws.Cells.Item(94, 13).Value = -35716073  # BSM!M94: -41668681 -> -35716073
$ws.Cells.Item(94, 14).Value = -4851.25  # BSM!N94: -4237 -> -4851.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1945.25  # BSM!H105: 1977.3549 -> 1945.25
$ws.Cells.Item(105, 9).Value = 1646.45  # BSM!I105: 1683.1052 -> 1646.45
$ws.Cells.Item(105, 11).Value = 1646.45  # BSM!K105: 1683.1052 -> 1646.45
$ws.Cells.Item(105, 13).Value = 100.55  # BSM!M105: 63.89480000000003 -> 100.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1440.12  # BSM!H107: 1562.6666 -> 1440.12
$ws.Cells.Item(107, 9).Value = 1421.9131  # BSM!I107: 1553.5264 -> 1421.9131
$ws.Cells.Item(107, 11).Value = 1421.9131  # BSM!K107: 1553.5264 -> 1421.9131
$ws.Cells.Item(107, 13).Value = 498.0869  # BSM!M107: 366.4736 -> 498.0869

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 1169.5  # CRP!H16: 1272 -> 1169.5
$ws.Cells.Item(16, 9).Value = 988.3333  # CRP!I16: 1053.7142 -> 988.3333
$ws.Cells.Item(16, 11).Value = 988.3333  # CRP!K16: 1053.7142 -> 988.3333
$ws.Cells.Item(16, 13).Value = -701.3333  # CRP!M16: -766.7141999999999 -> -701.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(113, 8).Value = 1169.5  # CRP!H113: 1272 -> 1169.5
$ws.Cells.Item(113, 9).Value = 988.3333  # CRP!I113: 1053.7142 -> 988.3333
$ws.Cells.Item(113, 11).Value = 988.3333  # CRP!K113: 1053.7142 -> 988.3333
$ws.Cells.Item(113, 13).Value = 1181.6667  # CRP!M113: 1116.2858 -> 1181.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(122, 8).Value = 2164.0908  # CRP!H122: 2638.375 -> 2164.0908
$ws.Cells.Item(122, 9).Value = 2765.5  # CRP!I122: 3698.5 -> 2765.5
$ws.Cells.Item(122, 10).Value = 1442.4  # CRP!J122: 1578.25 -> 1442.4
$ws.Cells.Item(122, 11).Value = 8296.5  # CRP!K122: 11095.5 -> 8296.5
$ws.Cells.Item(122, 12).Value = 4327.200000000001  # CRP!L122: 4734.75 -> 4327.200000000001
$ws.Cells.Item(122, 13).Value = -5846.5  # CRP!M122: -8645.5 -> -5846.5
$ws.Cells.Item(122, 14).Value = -9227.200000000001  # CRP!N122: -9634.75 -> -9227.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 1498.6666  # CRP!H132: 1591.5312 -> 1498.6666
$ws.Cells.Item(132, 9).Value = 1115.2858  # CRP!I132: 1175.2084 -> 1115.2858
$ws.Cells.Item(132, 11).Value = 3345.8574  # CRP!K132: 3525.6252 -> 3345.8574
$ws.Cells.Item(132, 13).Value = -815.8574000000003  # CRP!M132: -995.6251999999999 -> -815.8574000000003

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(134, 8).Value = 2810.8333  # CRP!H134: 3022.6667 -> 2810.8333
$ws.Cells.Item(134, 9).Value = 1579.9166  # CRP!I134: 1679.8182 -> 1579.9166
$ws.Cells.Item(134, 10).Value = 5272.6665  # CRP!J134: 5708.364 -> 5272.6665
$ws.Cells.Item(134, 11).Value = 4739.7498  # CRP!K134: 5039.4546 -> 4739.7498
$ws.Cells.Item(134, 12).Value = 15817.9995  # CRP!L134: 17125.092 -> 15817.9995
$ws.Cells.Item(134, 13).Value = -2204.7498  # CRP!M134: -2504.4546 -> -2204.7498
$ws.Cells.Item(134, 14).Value = -20887.9995  # CRP!N134: -22195.092 -> -20887.9995

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(135, 8).Value = 61805.285  # CRP!H135: 62912.57 -> 61805.285
$ws.Cells.Item(135, 10).Value = 61805.285  # CRP!J135: 62912.57 -> 61805.285
$ws.Cells.Item(135, 12).Value = 61805.285  # CRP!L135: 62912.57 -> 61805.285
$ws.Cells.Item(135, 14).Value = -71945.285  # CRP!N135: -73052.57000000001 -> -71945.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(117, 8).Value = 715  # CUL!H117: 731.8 -> 715
$ws.Cells.Item(117, 9).Value = 168.75  # CUL!I117: 165.44444 -> 168.75
$ws.Cells.Item(117, 10).Value = 2900  # CUL!J117: 1581.3334 -> 2900
$ws.Cells.Item(117, 11).Value = 506.25  # CUL!K117: 496.33332 -> 506.25
$ws.Cells.Item(117, 12).Value = 8700  # CUL!L117: 4744.0002 -> 8700
$ws.Cells.Item(117, 13).Value = 2935.75  # CUL!M117: 2945.66668 -> 2935.75
$ws.Cells.Item(117, 14).Value = -15584  # CUL!N117: -11628.0002 -> -15584

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 7161.7417  # GSM!H126: 7759.8887 -> 7161.7417
$ws.Cells.Item(126, 9).Value = 9574.764999999999  # GSM!I126: 11674.923 -> 9574.764999999999
$ws.Cells.Item(126, 10).Value = 4231.643  # GSM!J126: 4124.5 -> 4231.643
$ws.Cells.Item(126, 11).Value = 28724.295  # GSM!K126: 35024.769 -> 28724.295
$ws.Cells.Item(126, 12).Value = 12694.929  # GSM!L126: 12373.5 -> 12694.929
$ws.Cells.Item(126, 13).Value = -26254.295  # GSM!M126: -32554.769 -> -26254.295
$ws.Cells.Item(126, 14).Value = -17634.929  # GSM!N126: -17313.5 -> -17634.929

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2875.1  # GSM!H132: 3046.7036 -> 2875.1
$ws.Cells.Item(132, 9).Value = 2975.9524  # GSM!I132: 3250.1667 -> 2975.9524
$ws.Cells.Item(132, 11).Value = 8927.8572  # GSM!K132: 9750.500100000001 -> 8927.8572
$ws.Cells.Item(132, 13).Value = -6397.8572  # GSM!M132: -7220.500100000001 -> -6397.8572

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 61824.75  # GSM!H133: 64149.75 -> 61824.75
$ws.Cells.Item(133, 10).Value = 61824.75  # GSM!J133: 64149.75 -> 61824.75
$ws.Cells.Item(133, 12).Value = 61824.75  # GSM!L133: 64149.75 -> 61824.75
$ws.Cells.Item(133, 14).Value = -71944.75  # GSM!N133: -74269.75 -> -71944.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2145678  # LTW!H40: 2502991.5 -> 2145678
$ws.Cells.Item(40, 9).Value = 2502549.2  # LTW!I40: 3002699.8 -> 2502549.2
$ws.Cells.Item(40, 11).Value = 2502549.2  # LTW!K40: 3002699.8 -> 2502549.2
$ws.Cells.Item(40, 13).Value = -2502413.2  # LTW!M40: -3002563.8 -> -2502413.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(43, 8).Value = 2130000  # LTW!H43: 2625000 -> 2130000
$ws.Cells.Item(43, 10).Value = 2537500  # LTW!J43: 3333333.2 -> 2537500
$ws.Cells.Item(43, 12).Value = 2537500  # LTW!L43: 3333333.2 -> 2537500
$ws.Cells.Item(43, 14).Value = -2537886  # LTW!N43: -3333719.2 -> -2537886

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1993.0834  # LTW!H46: 1739.4762 -> 1993.0834
$ws.Cells.Item(46, 9).Value = 500.66666  # LTW!I46: 485.57144 -> 500.66666
$ws.Cells.Item(46, 10).Value = 2490.5557  # LTW!J46: 2366.4285 -> 2490.5557
$ws.Cells.Item(46, 11).Value = 500.66666  # LTW!K46: 485.57144 -> 500.66666
$ws.Cells.Item(46, 12).Value = 2490.5557  # LTW!L46: 2366.4285 -> 2490.5557
$ws.Cells.Item(46, 13).Value = -312.66666  # LTW!M46: -297.57144 -> -312.66666
$ws.Cells.Item(46, 14).Value = -2866.5557  # LTW!N46: -2742.4285 -> -2866.5557

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 3620.7778  # LTW!H61: 3949.125 -> 3620.7778
$ws.Cells.Item(61, 9).Value = 3948.5  # LTW!I61: 4370.5713 -> 3948.5
$ws.Cells.Item(61, 11).Value = 3948.5  # LTW!K61: 4370.5713 -> 3948.5
$ws.Cells.Item(61, 13).Value = -3746.5  # LTW!M61: -4168.5713 -> -3746.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(113, 8).Value = 3620.7778  # LTW!H113: 3949.125 -> 3620.7778
$ws.Cells.Item(113, 9).Value = 3948.5  # LTW!I113: 4370.5713 -> 3948.5
$ws.Cells.Item(113, 11).Value = 3948.5  # LTW!K113: 4370.5713 -> 3948.5
$ws.Cells.Item(113, 13).Value = -1778.5  # LTW!M113: -2200.5713 -> -1778.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 5396.3486  # LTW!H136: 5614.488 -> 5396.3486
$ws.Cells.Item(136, 9).Value = 5588.2856  # LTW!I136: 5825.2 -> 5588.2856
$ws.Cells.Item(136, 10).Value = 5213.136  # LTW!J136: 5413.8096 -> 5213.136
$ws.Cells.Item(136, 11).Value = 16764.8568  # LTW!K136: 17475.6 -> 16764.8568
$ws.Cells.Item(136, 12).Value = 15639.408  # LTW!L136: 16241.4288 -> 15639.408
$ws.Cells.Item(136, 13).Value = -14214.8568  # LTW!M136: -14925.6 -> -14214.8568
$ws.Cells.Item(136, 14).Value = -20739.408  # LTW!N136: -21341.4288 -> -20739.408

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 124250  # WVR!H46: 139666.67 -> 124250
$ws.Cells.Item(46, 10).Value = 124250  # WVR!J46: 139666.67 -> 124250
$ws.Cells.Item(46, 12).Value = 124250  # WVR!L46: 139666.67 -> 124250
$ws.Cells.Item(46, 14).Value = -124712  # WVR!N46: -140128.67 -> -124712

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value = 15700  # WVR!H98: 15645 -> 15700
$ws.Cells.Item(98, 10).Value = 15700  # WVR!J98: 15645 -> 15700
$ws.Cells.Item(98, 12).Value = 15700  # WVR!L98: 15645 -> 15700
$ws.Cells.Item(98, 14).Value = -21690  # WVR!N98: -21635 -> -21690

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(134, 8).Value = 124250  # WVR!H134: 139666.67 -> 124250
$ws.Cells.Item(134, 10).Value = 124250  # WVR!J134: 139666.67 -> 124250
$ws.Cells.Item(134, 12).Value = 372750  # WVR!L134: 419000.01 -> 372750
$ws.Cells.Item(134, 14).Value = -377820  # WVR!N134: -424070.01 -> -377820

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 7148.087  # WVR!H136: 7990.6 -> 7148.087
$ws.Cells.Item(136, 9).Value = 7948.1763  # WVR!I136: 9323.214 -> 7948.1763
$ws.Cells.Item(136, 11).Value = 23844.5289  # WVR!K136: 27969.642 -> 23844.5289
$ws.Cells.Item(136, 13).Value = -21294.5289  # WVR!M136: -25419.642 -> -21294.5289
